$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$ws.Cells.Item(2, 2).Value = 17.65992530191149
$ws.Cells.Item(3, 2).Value = 17.59826587873509
$ws.Cells.Item(4, 2).Value = 17.56418070193732
$ws.Cells.Item(5, 2).Value = 17.55124999241869
$ws.Cells.Item(6, 2).Value = 17.54916107488245
$ws.Cells.Item(7, 2).Value = 17.56400241704501
$ws.Cells.Item(8, 2).Value = 17.63788834326548
$ws.Cells.Item(9, 2).Value = 17.81222291131009
$ws.Cells.Item(10, 2).Value = 17.95751920730134
$ws.Cells.Item(11, 2).Value = 18.02717826812204
$ws.Cells.Item(12, 2).Value = 18.05405229901184
$ws.Cells.Item(13, 2).Value = 18.04824272580307
$ws.Cells.Item(14, 2).Value = 18.02937937571015
$ws.Cells.Item(15, 2).Value = 18.01788906843487
$ws.Cells.Item(16, 2).Value = 17.95303711237774
$ws.Cells.Item(17, 2).Value = 17.91415349985963
$ws.Cells.Item(18, 2).Value = 17.89212510014787
$ws.Cells.Item(19, 2).Value = 17.88472493926403
$ws.Cells.Item(20, 2).Value = 17.91825802645606
$ws.Cells.Item(21, 2).Value = 18.03490668408837
$ws.Cells.Item(22, 2).Value = 18.11402408170125
$ws.Cells.Item(23, 2).Value = 18.07153977393441
$ws.Cells.Item(24, 2).Value = 17.91640135231762
$ws.Cells.Item(25, 2).Value = 17.76198066449103

# Column D
$ws.Cells.Item(2, 4).Value = 8.020924470889165
$ws.Cells.Item(3, 4).Value = 8.036462109335512
$ws.Cells.Item(4, 4).Value = 8.046428708180949
$ws.Cells.Item(5, 4).Value = 8.050597692508742
$ws.Cells.Item(6, 4).Value = 8.0512964512177
$ws.Cells.Item(7, 4).Value = 8.04648449688119
$ws.Cells.Item(8, 4).Value = 8.026193516799053
$ws.Cells.Item(9, 4).Value = 7.989773115788501
$ws.Cells.Item(10, 4).Value = 7.965050479929932
$ws.Cells.Item(11, 4).Value = 7.95424144119725
$ws.Cells.Item(12, 4).Value = 7.950210949618322
$ws.Cells.Item(13, 4).Value = 7.951076205223134
$ws.Cells.Item(14, 4).Value = 7.953908595867916
$ws.Cells.Item(15, 4).Value = 7.955651669564171
$ws.Cells.Item(16, 4).Value = 7.96576565704294
$ws.Cells.Item(17, 4).Value = 7.9720821184684
$ws.Cells.Item(18, 4).Value = 7.975756365559164
$ws.Cells.Item(19, 4).Value = 7.977007483602813
$ws.Cells.Item(20, 4).Value = 7.971405459631488
$ws.Cells.Item(21, 4).Value = 7.953074955498085
$ws.Cells.Item(22, 4).Value = 7.941460013968896
$ws.Cells.Item(23, 4).Value = 7.947625798858011
$ws.Cells.Item(24, 4).Value = 7.971711243648715
$ws.Cells.Item(25, 4).Value = 7.999266992548979

# Column E
$ws.Cells.Item(2, 5).Value = 13.11632871265028
$ws.Cells.Item(3, 5).Value = 13.05772177745514
$ws.Cells.Item(4, 5).Value = 13.02076690090516
$ws.Cells.Item(5, 5).Value = 13.00546386829334
$ws.Cells.Item(6, 5).Value = 13.00290810559118
$ws.Cells.Item(7, 5).Value = 13.02056150482379
$ws.Cells.Item(8, 5).Value = 13.09632252678858
$ws.Cells.Item(9, 5).Value = 13.23718465553841
$ws.Cells.Item(10, 5).Value = 13.33597357710753
$ws.Cells.Item(11, 5).Value = 13.37989322500058
$ws.Cells.Item(12, 5).Value = 13.39637752454941
$ws.Cells.Item(13, 5).Value = 13.39283388604384
$ws.Cells.Item(14, 5).Value = 13.38125234929946
$ws.Cells.Item(15, 5).Value = 13.37413915109363
$ws.Cells.Item(16, 5).Value = 13.3330827625076
$ws.Cells.Item(17, 5).Value = 13.30763444742165
$ws.Cells.Item(18, 5).Value = 13.29290127640554
$ws.Cells.Item(19, 5).Value = 13.28789640820069
$ws.Cells.Item(20, 5).Value = 13.31035340003726
$ws.Cells.Item(21, 5).Value = 13.38465812373666
$ws.Cells.Item(22, 5).Value = 13.43236283072609
$ws.Cells.Item(23, 5).Value = 13.40698056235513
$ws.Cells.Item(24, 5).Value = 13.30912448115117
$ws.Cells.Item(25, 5).Value = 13.19990310304972

# Column F
$ws.Cells.Item(2, 6).Value = 40.8134782372903
$ws.Cells.Item(3, 6).Value = 40.59733323531989
$ws.Cells.Item(4, 6).Value = 40.47357619555765
$ws.Cells.Item(5, 6).Value = 40.42542975013981
$ws.Cells.Item(6, 6).Value = 40.41757401426891
$ws.Cells.Item(7, 6).Value = 40.47291758048351
$ws.Cells.Item(8, 6).Value = 40.73711345699579
$ws.Cells.Item(9, 6).Value = 41.32464978319835
$ws.Cells.Item(10, 6).Value = 41.7962677144971
$ws.Cells.Item(11, 6).Value = 42.0189259780097
$ws.Cells.Item(12, 6).Value = 42.10435611849238
$ws.Cells.Item(13, 6).Value = 42.08590851658845
$ws.Cells.Item(14, 6).Value = 42.02593234595955
$ws.Cells.Item(15, 6).Value = 41.98933871487876
$ws.Cells.Item(16, 6).Value = 41.78187533444379
$ws.Cells.Item(17, 6).Value = 41.65664528215967
$ws.Cells.Item(18, 6).Value = 41.58538358947834
$ws.Cells.Item(19, 6).Value = 41.56138896783393
$ws.Cells.Item(20, 6).Value = 41.66989719667355
$ws.Cells.Item(21, 6).Value = 42.04351899219343
$ws.Cells.Item(22, 6).Value = 42.29417133373933
$ws.Cells.Item(23, 6).Value = 42.15981958905218
$ws.Cells.Item(24, 6).Value = 41.66390371344575
$ws.Cells.Item(25, 6).Value = 41.15850447990817

# Column G
$ws.Cells.Item(2, 7).Value = 3.718827194612049
$ws.Cells.Item(3, 7).Value = 3.722364706182779
$ws.Cells.Item(4, 7).Value = 3.724649971851191
$ws.Cells.Item(5, 7).Value = 3.725609810881213
$ws.Cells.Item(6, 7).Value = 3.725770920342412
$ws.Cells.Item(7, 7).Value = 3.724662800740326
$ws.Cells.Item(8, 7).Value = 3.72002349315398
$ws.Cells.Item(9, 7).Value = 3.711819389360516
$ws.Cells.Item(10, 7).Value = 3.706329887565656
$ws.Cells.Item(11, 7).Value = 3.703947977059321
$ws.Cells.Item(12, 7).Value = 3.703062478895439
$ws.Cells.Item(13, 7).Value = 3.703252455404783
$ws.Cells.Item(14, 7).Value = 3.70387479688213
$ws.Cells.Item(15, 7).Value = 3.704258142209803
$ws.Cells.Item(16, 7).Value = 3.70648786249335
$ws.Cells.Item(17, 7).Value = 3.707885181879818
$ws.Cells.Item(18, 7).Value = 3.708699740701884
$ws.Cells.Item(19, 7).Value = 3.708977404164868
$ws.Cells.Item(20, 7).Value = 3.707735311841196
$ws.Cells.Item(21, 7).Value = 3.703691553638604
$ws.Cells.Item(22, 7).Value = 3.701144734257809
$ws.Cells.Item(23, 7).Value = 3.702495267356639
$ws.Cells.Item(24, 7).Value = 3.707803033124173
$ws.Cells.Item(25, 7).Value = 3.71394384776582

# Column K
$ws.Cells.Item(2, 11).Value = 14.16554261532486
$ws.Cells.Item(3, 11).Value = 13.91046660048228
$ws.Cells.Item(4, 11).Value = 13.75605059585088
$ws.Cells.Item(5, 11).Value = 13.69377099134552
$ws.Cells.Item(6, 11).Value = 13.68347123288239
$ws.Cells.Item(7, 11).Value = 13.75520793309046
$ws.Cells.Item(8, 11).Value = 14.07719102307844
$ws.Cells.Item(9, 11).Value = 14.72200927057147
$ws.Cells.Item(10, 11).Value = 15.19857555524722
$ws.Cells.Item(11, 11).Value = 15.4148981223927
$ws.Cells.Item(12, 11).Value = 15.49665961166916
$ws.Cells.Item(13, 11).Value = 15.47905906631434
$ws.Cells.Item(14, 11).Value = 15.42162828702746
$ws.Cells.Item(15, 11).Value = 15.38642753520474
$ws.Cells.Item(16, 11).Value = 15.18442133513556
$ws.Cells.Item(17, 11).Value = 15.0603131945544
$ws.Cells.Item(18, 11).Value = 14.98889026396233
$ws.Cells.Item(19, 11).Value = 14.96470368478576
$ws.Cells.Item(20, 11).Value = 15.07352942318822
$ws.Cells.Item(21, 11).Value = 15.43850198905788
$ws.Cells.Item(22, 11).Value = 15.67609041991278
$ws.Cells.Item(23, 11).Value = 15.54939912614882
$ws.Cells.Item(24, 11).Value = 15.06755458223241
$ws.Cells.Item(25, 11).Value = 14.54669363785277

# Column L
$ws.Cells.Item(2, 12).Value = 9.831300910267329
$ws.Cells.Item(3, 12).Value = 9.824048139237986
$ws.Cells.Item(4, 12).Value = 9.821481016222485
$ws.Cells.Item(5, 12).Value = 9.82091025501113
$ws.Cells.Item(6, 12).Value = 9.82084421762775
$ws.Cells.Item(7, 12).Value = 9.821471392796472
$ws.Cells.Item(8, 12).Value = 9.828409421914671
$ws.Cells.Item(9, 12).Value = 9.856922232985941
$ws.Cells.Item(10, 12).Value = 9.886865478190778
$ws.Cells.Item(11, 12).Value = 9.90241506014805
$ws.Cells.Item(12, 12).Value = 9.90857790755404
$ws.Cells.Item(13, 12).Value = 9.907238463337041
$ws.Cells.Item(14, 12).Value = 9.902916595335851
$ws.Cells.Item(15, 12).Value = 9.90030499644914
$ws.Cells.Item(16, 12).Value = 9.885887851460401
$ws.Cells.Item(17, 12).Value = 9.877535423212629
$ws.Cells.Item(18, 12).Value = 9.872913004179077
$ws.Cells.Item(19, 12).Value = 9.871379214013674
$ws.Cells.Item(20, 12).Value = 9.878405769470586
$ws.Cells.Item(21, 12).Value = 9.904178604368658
$ws.Cells.Item(22, 12).Value = 9.922621519808825
$ws.Cells.Item(23, 12).Value = 9.912632861921473
$ws.Cells.Item(24, 12).Value = 9.878011726592888
$ws.Cells.Item(25, 12).Value = 9.847621308953348

# Column M
$ws.Cells.Item(2, 13).Value = 16.58706324801927
$ws.Cells.Item(3, 13).Value = 16.58428121181482
$ws.Cells.Item(4, 13).Value = 16.58550094974628
$ws.Cells.Item(5, 13).Value = 16.58673559017334
$ws.Cells.Item(6, 13).Value = 16.5869851825535
$ws.Cells.Item(7, 13).Value = 16.58551461256742
$ws.Cells.Item(8, 13).Value = 16.58549720519202
$ws.Cells.Item(9, 13).Value = 16.6086246367455
$ws.Cells.Item(10, 13).Value = 16.63962538393956
$ws.Cells.Item(11, 13).Value = 16.656739944318
$ws.Cells.Item(12, 13).Value = 16.66365073410724
$ws.Cells.Item(13, 13).Value = 16.66214330865871
$ws.Cells.Item(14, 13).Value = 16.65729990371518
$ws.Cells.Item(15, 13).Value = 16.65438905817495
$ws.Cells.Item(16, 13).Value = 16.63856724808725
$ws.Cells.Item(17, 13).Value = 16.62963035159
$ws.Cells.Item(18, 13).Value = 16.6247738219082
$ws.Cells.Item(19, 13).Value = 16.62317831039424
$ws.Cells.Item(20, 13).Value = 16.63055235644277
$ws.Cells.Item(21, 13).Value = 16.65871089015656
$ws.Cells.Item(22, 13).Value = 16.67961801572945
$ws.Cells.Item(23, 13).Value = 16.66823155084776
$ws.Cells.Item(24, 13).Value = 16.6301346412416
$ws.Cells.Item(25, 13).Value = 16.59989854763121

# Column N
$ws.Cells.Item(2, 14).Value = 24.13133141334478
$ws.Cells.Item(3, 14).Value = 24.1671897247264
$ws.Cells.Item(4, 14).Value = 24.19104768919717
$ws.Cells.Item(5, 14).Value = 24.2012327799564
$ws.Cells.Item(6, 14).Value = 24.20295195389826
$ws.Cells.Item(7, 14).Value = 24.19118317533733
$ws.Cells.Item(8, 14).Value = 24.14331322984315
$ws.Cells.Item(9, 14).Value = 24.06405183769371
$ws.Cells.Item(10, 14).Value = 24.01473394500875
$ws.Cells.Item(11, 14).Value = 23.99423577991284
$ws.Cells.Item(12, 14).Value = 23.98675236385237
$ws.Cells.Item(13, 14).Value = 23.98835164731242
$ws.Cells.Item(14, 14).Value = 23.99361452566237
$ws.Cells.Item(15, 14).Value = 23.9968745064322
$ws.Cells.Item(16, 14).Value = 24.01611254579517
$ws.Cells.Item(17, 14).Value = 24.02841071994871
$ws.Cells.Item(18, 14).Value = 24.03566658700933
$ws.Cells.Item(19, 14).Value = 24.03815459748149
$ws.Cells.Item(20, 14).Value = 24.02708269083336
$ws.Cells.Item(21, 14).Value = 23.9920611220087
$ws.Cells.Item(22, 14).Value = 23.97079766535857
$ws.Cells.Item(23, 14).Value = 23.98199757365515
$ws.Cells.Item(24, 14).Value = 24.02768251501747
$ws.Cells.Item(25, 14).Value = 24.08392909709781

